$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1248.7778  # H15: 2490.2715 -> 1248.7778
$ws.Cells.Item(15, 9).Value = 1248.7778  # I15: 2490.2715 -> 1248.7778
$ws.Cells.Item(15, 11).Value = 3746.3334  # K15: 7470.814499999999 -> 3746.3334
$ws.Cells.Item(15, 13).Value = -3577.3334  # M15: -7301.814499999999 -> -3577.3334

$ws.Cells.Item(32, 8).Value = 1099.3334  # H32: 1234.2858 -> 1099.3334
$ws.Cells.Item(32, 9).Value = 823  # I32: 879.3333 -> 823
$ws.Cells.Item(32, 10).Value = 1320.4  # J32: 1500.5 -> 1320.4
$ws.Cells.Item(32, 11).Value = 823  # K32: 879.3333 -> 823
$ws.Cells.Item(32, 12).Value = 1320.4  # L32: 1500.5 -> 1320.4
$ws.Cells.Item(32, 13).Value = -497  # M32: -553.3333 -> -497
$ws.Cells.Item(32, 14).Value = -1972.4  # N32: -2152.5 -> -1972.4

$ws.Cells.Item(86, 8).Value = 2714.1428  # H86: 3499.75 -> 2714.1428
$ws.Cells.Item(86, 9).Value = 2375  # I86: 3000 -> 2375
$ws.Cells.Item(86, 10).Value = 3166.3333  # J86: 3999.5 -> 3166.3333
$ws.Cells.Item(86, 11).Value = 2375  # K86: 3000 -> 2375
$ws.Cells.Item(86, 12).Value = 3166.3333  # L86: 3999.5 -> 3166.3333
$ws.Cells.Item(86, 13).Value = -1252  # M86: -1877 -> -1252
$ws.Cells.Item(86, 14).Value = -5412.3333  # N86: -6245.5 -> -5412.3333

$ws.Cells.Item(89, 8).Value = 2714.1428  # H89: 3499.75 -> 2714.1428
$ws.Cells.Item(89, 9).Value = 2375  # I89: 3000 -> 2375
$ws.Cells.Item(89, 10).Value = 3166.3333  # J89: 3999.5 -> 3166.3333
$ws.Cells.Item(89, 11).Value = 11875  # K89: 15000 -> 11875
$ws.Cells.Item(89, 12).Value = 15831.6665  # L89: 19997.5 -> 15831.6665
$ws.Cells.Item(89, 13).Value = -6259  # M89: -9384 -> -6259
$ws.Cells.Item(89, 14).Value = -27063.6665  # N89: -31229.5 -> -27063.6665

$ws.Cells.Item(123, 8).Value = 66024.96000000001  # H123: 68700 -> 66024.96000000001
$ws.Cells.Item(123, 10).Value = 66024.96000000001  # J123: 68700 -> 66024.96000000001
$ws.Cells.Item(123, 12).Value = 66024.96000000001  # L123: 68700 -> 66024.96000000001
$ws.Cells.Item(123, 14).Value = -75824.96000000001  # N123: -78500 -> -75824.96000000001

$ws.Cells.Item(127, 8).Value = 937.24  # H127: 929.78 -> 937.24
$ws.Cells.Item(127, 9).Value = 331.8889  # I127: 323.9 -> 331.8889
$ws.Cells.Item(127, 10).Value = 997.10986  # J127: 997.1 -> 997.10986
$ws.Cells.Item(127, 11).Value = 995.6667  # K127: 971.6999999999999 -> 995.6667
$ws.Cells.Item(127, 12).Value = 2991.32958  # L127: 2991.3 -> 2991.32958
$ws.Cells.Item(127, 13).Value = 3964.3333  # M127: 3988.3 -> 3964.3333
$ws.Cells.Item(127, 14).Value = -12911.32958  # N127: -12911.3 -> -12911.32958

$ws.Cells.Item(129, 8).Value = 936.78687  # H129: 944.35 -> 936.78687
$ws.Cells.Item(129, 9).Value = 498.4  # I129: 491.85715 -> 498.4
$ws.Cells.Item(129, 10).Value = 1079.7391  # J129: 1188 -> 1079.7391
$ws.Cells.Item(129, 11).Value = 1495.2  # K129: 1475.57145 -> 1495.2
$ws.Cells.Item(129, 12).Value = 3239.2173  # L129: 3564 -> 3239.2173
$ws.Cells.Item(129, 13).Value = 3504.8  # M129: 3524.42855 -> 3504.8
$ws.Cells.Item(129, 14).Value = -13239.2173  # N129: -13564 -> -13239.2173

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 28735.143  # H32: 12100.923 -> 28735.143
$ws.Cells.Item(32, 9).Value = 29921.9  # I32: 13128.787 -> 29921.9
$ws.Cells.Item(32, 10).Value = 5000  # J32: 2439 -> 5000
$ws.Cells.Item(32, 11).Value = 29921.9  # K32: 13128.787 -> 29921.9
$ws.Cells.Item(32, 12).Value = 5000  # L32: 2439 -> 5000
$ws.Cells.Item(32, 13).Value = -29634.9  # M32: -12841.787 -> -29634.9
$ws.Cells.Item(32, 14).Value = -5574  # N32: -3013 -> -5574

$ws.Cells.Item(63, 8).Value = 4091  # H63: 4365 -> 4091
$ws.Cells.Item(63, 9).Value = 4742.857  # I63: 4755 -> 4742.857
$ws.Cells.Item(63, 10).Value = 2950.25  # J63: 3000 -> 2950.25
$ws.Cells.Item(63, 11).Value = 4742.857  # K63: 4755 -> 4742.857
$ws.Cells.Item(63, 12).Value = 2950.25  # L63: 3000 -> 2950.25
$ws.Cells.Item(63, 13).Value = -4056.857  # M63: -4069 -> -4056.857
$ws.Cells.Item(63, 14).Value = -4322.25  # N63: -4372 -> -4322.25

$ws.Cells.Item(66, 8).Value = 4091  # H66: 4365 -> 4091
$ws.Cells.Item(66, 9).Value = 4742.857  # I66: 4755 -> 4742.857
$ws.Cells.Item(66, 10).Value = 2950.25  # J66: 3000 -> 2950.25
$ws.Cells.Item(66, 11).Value = 23714.285  # K66: 23775 -> 23714.285
$ws.Cells.Item(66, 12).Value = 14751.25  # L66: 15000 -> 14751.25
$ws.Cells.Item(66, 13).Value = -20282.285  # M66: -20343 -> -20282.285
$ws.Cells.Item(66, 14).Value = -21615.25  # N66: -21864 -> -21615.25

$ws.Cells.Item(74, 8).Value = 1103.3226  # H74: 1113.4333 -> 1103.3226
$ws.Cells.Item(74, 9).Value = 888.13635  # I74: 894.9545000000001 -> 888.13635
$ws.Cells.Item(74, 10).Value = 1629.3334  # J74: 1714.25 -> 1629.3334
$ws.Cells.Item(74, 11).Value = 888.13635  # K74: 894.9545000000001 -> 888.13635
$ws.Cells.Item(74, 12).Value = 1629.3334  # L74: 1714.25 -> 1629.3334
$ws.Cells.Item(74, 13).Value = -14.13634999999999  # M74: -20.95450000000005 -> -14.13634999999999
$ws.Cells.Item(74, 14).Value = -3377.3334  # N74: -3462.25 -> -3377.3334

$ws.Cells.Item(77, 8).Value = 1103.3226  # H77: 1113.4333 -> 1103.3226
$ws.Cells.Item(77, 9).Value = 888.13635  # I77: 894.9545000000001 -> 888.13635
$ws.Cells.Item(77, 10).Value = 1629.3334  # J77: 1714.25 -> 1629.3334
$ws.Cells.Item(77, 11).Value = 4440.68175  # K77: 4474.7725 -> 4440.68175
$ws.Cells.Item(77, 12).Value = 8146.666999999999  # L77: 8571.25 -> 8146.666999999999
$ws.Cells.Item(77, 13).Value = -72.68174999999974  # M77: -106.7725 -> -72.68174999999974
$ws.Cells.Item(77, 14).Value = -16882.667  # N77: -17307.25 -> -16882.667

$ws.Cells.Item(132, 8).Value = 3024.8125  # H132: 2616.8235 -> 3024.8125
$ws.Cells.Item(132, 9).Value = 2112.5  # I132: 1600.4 -> 2112.5
$ws.Cells.Item(132, 10).Value = 3937.125  # J132: 4068.8572 -> 3937.125
$ws.Cells.Item(132, 11).Value = 6337.5  # K132: 4801.200000000001 -> 6337.5
$ws.Cells.Item(132, 12).Value = 11811.375  # L132: 12206.5716 -> 11811.375
$ws.Cells.Item(132, 13).Value = -3807.5  # M132: -2271.200000000001 -> -3807.5
$ws.Cells.Item(132, 14).Value = -16871.375  # N132: -17266.5716 -> -16871.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 10000  # H15: 0 -> 10000
$ws.Cells.Item(15, 10).Value = 10000  # J15: 0 -> 10000
$ws.Cells.Item(15, 12).Value = 10000  # L15: 0 -> 10000
$ws.Cells.Item(15, 14).Value = -10454  # N15: None -> -10454

$ws.Cells.Item(55, 8).Value = 24500  # H55: 0 -> 24500
$ws.Cells.Item(55, 10).Value = 24500  # J55: 0 -> 24500
$ws.Cells.Item(55, 12).Value = 24500  # L55: 0 -> 24500
$ws.Cells.Item(55, 14).Value = -25046  # N55: None -> -25046

$ws.Cells.Item(86, 8).Value = 168622.5  # H86: 84450.336 -> 168622.5
$ws.Cells.Item(86, 9).Value = 2345.6  # I86: 1249.7 -> 2345.6
$ws.Cells.Item(86, 10).Value = 1000007  # J86: 500453.5 -> 1000007
$ws.Cells.Item(86, 11).Value = 2345.6  # K86: 1249.7 -> 2345.6
$ws.Cells.Item(86, 12).Value = 1000007  # L86: 500453.5 -> 1000007
$ws.Cells.Item(86, 13).Value = -1222.6  # M86: -126.7 -> -1222.6
$ws.Cells.Item(86, 14).Value = -1002253  # N86: -502699.5 -> -1002253

$ws.Cells.Item(89, 8).Value = 168622.5  # H89: 84450.336 -> 168622.5
$ws.Cells.Item(89, 9).Value = 2345.6  # I89: 1249.7 -> 2345.6
$ws.Cells.Item(89, 10).Value = 1000007  # J89: 500453.5 -> 1000007
$ws.Cells.Item(89, 11).Value = 11728  # K89: 6248.5 -> 11728
$ws.Cells.Item(89, 12).Value = 5000035  # L89: 2502267.5 -> 5000035
$ws.Cells.Item(89, 13).Value = -6112  # M89: -632.5 -> -6112
$ws.Cells.Item(89, 14).Value = -5011267  # N89: -2513499.5 -> -5011267

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(106, 8).Value = 152225  # H106: 500000 -> 152225
$ws.Cells.Item(106, 10).Value = 152225  # J106: 500000 -> 152225
$ws.Cells.Item(106, 12).Value = 152225  # L106: 500000 -> 152225
$ws.Cells.Item(106, 14).Value = -154749  # N106: -502524 -> -154749

$ws.Cells.Item(118, 8).Value = 211111  # H118: 117805.5 -> 211111
$ws.Cells.Item(118, 10).Value = 211111  # J118: 117805.5 -> 211111
$ws.Cells.Item(118, 12).Value = 211111  # L118: 117805.5 -> 211111
$ws.Cells.Item(118, 14).Value = -214425  # N118: -121119.5 -> -214425

$ws.Cells.Item(122, 8).Value = 2035.5714  # H122: 2524.3794 -> 2035.5714
$ws.Cells.Item(122, 9).Value = 1942.44  # I122: 2460.3333 -> 1942.44
$ws.Cells.Item(122, 10).Value = 2268.4  # J122: 2692.5 -> 2268.4
$ws.Cells.Item(122, 11).Value = 5827.32  # K122: 7380.999899999999 -> 5827.32
$ws.Cells.Item(122, 12).Value = 6805.200000000001  # L122: 8077.5 -> 6805.200000000001
$ws.Cells.Item(122, 13).Value = -3377.32  # M122: -4930.999899999999 -> -3377.32
$ws.Cells.Item(122, 14).Value = -11705.2  # N122: -12977.5 -> -11705.2

$ws.Cells.Item(132, 8).Value = 1980.6666  # H132: 1490.125 -> 1980.6666
$ws.Cells.Item(132, 9).Value = 1642.5834  # I132: 1337.5652 -> 1642.5834
$ws.Cells.Item(132, 10).Value = 3333  # J132: 4999 -> 3333
$ws.Cells.Item(132, 11).Value = 4927.7502  # K132: 4012.6956 -> 4927.7502
$ws.Cells.Item(132, 12).Value = 9999  # L132: 14997 -> 9999
$ws.Cells.Item(132, 13).Value = -2397.7502  # M132: -1482.6956 -> -2397.7502
$ws.Cells.Item(132, 14).Value = -15059  # N132: -20057 -> -15059

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 4808647  # H2: 1494.5 -> 4808647
$ws.Cells.Item(2, 9).Value = 1515.7142  # I2: 2152 -> 1515.7142
$ws.Cells.Item(2, 10).Value = 10416966  # J2: 398.66666 -> 10416966
$ws.Cells.Item(2, 11).Value = 9094.285199999998  # K2: 12912 -> 9094.285199999998
$ws.Cells.Item(2, 12).Value = 62501796  # L2: 2391.99996 -> 62501796
$ws.Cells.Item(2, 13).Value = -8981.285199999998  # M2: -12799 -> -8981.285199999998
$ws.Cells.Item(2, 14).Value = -62502022  # N2: -2617.99996 -> -62502022

$ws.Cells.Item(12, 8).Value = 743265.4  # H12: 772995.6 -> 743265.4
$ws.Cells.Item(12, 9).Value = 90.2  # I12: 110.25 -> 90.2
$ws.Cells.Item(12, 11).Value = 270.6  # K12: 330.75 -> 270.6
$ws.Cells.Item(12, 13).Value = -97.60000000000002  # M12: -157.75 -> -97.60000000000002

$ws.Cells.Item(23, 8).Value = 384.64285  # H23: 391.23077 -> 384.64285
$ws.Cells.Item(23, 10).Value = 398.75  # J23: 407.81818 -> 398.75
$ws.Cells.Item(23, 12).Value = 1196.25  # L23: 1223.45454 -> 1196.25
$ws.Cells.Item(23, 14).Value = -1666.25  # N23: -1693.45454 -> -1666.25

$ws.Cells.Item(88, 8).Value = 85000  # H88: 140000 -> 85000
$ws.Cells.Item(88, 10).Value = 85000  # J88: 140000 -> 85000
$ws.Cells.Item(88, 12).Value = 255000  # L88: 420000 -> 255000
$ws.Cells.Item(88, 14).Value = -255856  # N88: -420856 -> -255856

$ws.Cells.Item(91, 8).Value = 85000  # H91: 140000 -> 85000
$ws.Cells.Item(91, 10).Value = 85000  # J91: 140000 -> 85000
$ws.Cells.Item(91, 12).Value = 255000  # L91: 420000 -> 255000
$ws.Cells.Item(91, 14).Value = -257964  # N91: -422964 -> -257964

$ws.Cells.Item(94, 8).Value = 3743.6667  # H94: 3667.125 -> 3743.6667
$ws.Cells.Item(94, 9).Value = 0  # I94: 2605 -> 0
$ws.Cells.Item(94, 10).Value = 3743.6667  # J94: 4021.1667 -> 3743.6667
$ws.Cells.Item(94, 11).Value = 0  # K94: 7815 -> 0
$ws.Cells.Item(94, 12).Value = 11231.0001  # L94: 12063.5001 -> 11231.0001
$ws.Cells.Item(94, 13).ClearContents()  # M94: -7139 -> (removed)
$ws.Cells.Item(94, 14).Value = -12583.0001  # N94: -13415.5001 -> -12583.0001

$ws.Cells.Item(99, 8).Value = 3533.3333  # H99: 5481.25 -> 3533.3333
$ws.Cells.Item(99, 9).Value = 300  # I99: 962.5 -> 300
$ws.Cells.Item(99, 11).Value = 900  # K99: 2887.5 -> 900
$ws.Cells.Item(99, 13).Value = 1346  # M99: -641.5 -> 1346

$ws.Cells.Item(106, 8).Value = 8333  # H106: 6666.6665 -> 8333
$ws.Cells.Item(106, 10).Value = 8333  # J106: 6666.6665 -> 8333
$ws.Cells.Item(106, 12).Value = 24999  # L106: 19999.9995 -> 24999
$ws.Cells.Item(106, 14).Value = -26891  # N106: -21891.9995 -> -26891

$ws.Cells.Item(107, 8).Value = 659.0909  # H107: 660 -> 659.0909
$ws.Cells.Item(107, 10).Value = 371.66666  # J107: 373.33334 -> 371.66666
$ws.Cells.Item(107, 12).Value = 1114.99998  # L107: 1120.00002 -> 1114.99998
$ws.Cells.Item(107, 14).Value = -4954.999980000001  # N107: -4960.000019999999 -> -4954.999980000001

$ws.Cells.Item(131, 8).Value = 19232374  # H131: 12347034 -> 19232374
$ws.Cells.Item(131, 10).Value = 19232374  # J131: 12347034 -> 19232374
$ws.Cells.Item(131, 12).Value = 57697122  # L131: 37041102 -> 57697122
$ws.Cells.Item(131, 14).Value = -57707202  # N131: -37051182 -> -57707202

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0  # H5: 9800 -> 0
$ws.Cells.Item(5, 9).Value = 0  # I5: 9800 -> 0
$ws.Cells.Item(5, 11).Value = 0  # K5: 9800 -> 0
$ws.Cells.Item(5, 13).ClearContents()  # M5: -9688 -> (removed)

$ws.Cells.Item(105, 8).Value = 0  # H105: 48450 -> 0
$ws.Cells.Item(105, 10).Value = 0  # J105: 48450 -> 0
$ws.Cells.Item(105, 12).Value = 0  # L105: 48450 -> 0
$ws.Cells.Item(105, 14).ClearContents()  # N105: -55438 -> (removed)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4497.5  # H62: 4553.3335 -> 4497.5

$ws.Cells.Item(65, 8).Value = 4497.5  # H65: 4553.3335 -> 4497.5

$ws.Cells.Item(137, 8).Value = 21166.666  # H137: 34993.332 -> 21166.666
$ws.Cells.Item(137, 10).Value = 21166.666  # J137: 34993.332 -> 21166.666
$ws.Cells.Item(137, 12).Value = 21166.666  # L137: 34993.332 -> 21166.666
$ws.Cells.Item(137, 14).Value = -31366.666  # N137: -45193.332 -> -31366.666
